# Apply the updated crypto price/volume snapshot (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.365.18"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "2.408.58"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'506.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.07%  "
$ws.Range("D6").Value = "'133.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("D8").Value = "'0.557"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").Value = "2.443.09"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").Value = "'0.0977"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("D12").Value = "'0.322"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("E13").Value = "  -7.24%  "
$ws.Range("D14").Value = "2.847.93"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "57.199.84"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "'21.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "2.468.57"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "'10.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "'4.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "'313.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "'6.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.91%  "
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").Value = "'65.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").Value = "'0.993"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").Value = "2.533.01"
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("D28").Value = "'0.381"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.37%  "
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("D30").Value = "'7.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.15%  "
$ws.Range("D31").Value = "'173.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "0.0₃0728"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D37").Value = "'0.991"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("E39").Value = "  +3.77%  "
$ws.Range("D40").Value = "'3.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").Value = "'36.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "'0.812"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "'134.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.93%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'5.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.06%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'3.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").Value = "'256.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("E48").Value = "  -2.78%  "
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("E51").Value = "  +0.71%  "
